# Update "Datos crudos" sheet with the new temperature measurement run.
# The original 32 rows of data (rows 2-33) are replaced with 37 new rows
# (rows 2-38) taken from a later measurement session (2023-12-12), and the
# "Fila final" input cell (H2) / formula in H3 are adjusted to match.

$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("Datos crudos")
$wsValid = $wb.Worksheets.Item("Datos válidos")

# New timestamp strings for column C (rows 2..38)
$timestamps = @(
    "2023-12-12 03:20:30",
    "2023-12-12 03:21:30",
    "2023-12-12 03:22:31",
    "2023-12-12 03:23:32",
    "2023-12-12 03:24:33",
    "2023-12-12 03:25:34",
    "2023-12-12 03:26:35",
    "2023-12-12 03:27:36",
    "2023-12-12 03:28:37",
    "2023-12-12 03:29:38",
    "2023-12-12 03:30:39",
    "2023-12-12 03:31:40",
    "2023-12-12 03:32:41",
    "2023-12-12 03:33:42",
    "2023-12-12 03:34:43",
    "2023-12-12 03:35:44",
    "2023-12-12 03:36:45",
    "2023-12-12 03:37:46",
    "2023-12-12 03:38:47",
    "2023-12-12 03:39:48",
    "2023-12-12 03:40:49",
    "2023-12-12 03:41:50",
    "2023-12-12 03:42:51",
    "2023-12-12 03:43:52",
    "2023-12-12 03:44:53",
    "2023-12-12 03:45:54",
    "2023-12-12 03:46:55",
    "2023-12-12 03:47:56",
    "2023-12-12 03:48:57",
    "2023-12-12 03:49:58",
    "2023-12-12 03:50:59",
    "2023-12-12 03:52:00",
    "2023-12-12 03:53:01",
    "2023-12-12 03:54:02",
    "2023-12-12 03:55:03",
    "2023-12-12 03:56:04",
    "2023-12-12 03:57:04"
)

# New temperature readings for column E (rows 2..38)
$temperatures = @(
    27.377622377622401,
    26.722027972027998,
    26.1975524475524,
    26.1975524475524,
    25.8041958041958,
    25.541958041958001,
    25.1486013986014,
    25.1486013986014,
    25.1486013986014,
    24.7552447552447,
    24.7552447552447,
    23.575174825174798,
    24.7552447552447,
    24.7552447552447,
    24.624125874125799,
    24.361888111888099,
    24.624125874125799,
    23.837412587412501,
    23.181818181818102,
    23.312937062936999,
    23.312937062936999,
    23.575174825174798,
    23.181818181818102,
    23.0506993006993,
    24.0996503496503,
    23.181818181818102,
    23.444055944055901,
    23.0506993006993,
    23.312937062936999,
    23.575174825174798,
    24.361888111888099,
    23.181818181818102,
    24.230769230769202,
    25.410839160839199,
    25.8041958041958,
    26.1975524475524,
    26.328671328671302
)

$startRow = 2
$count = $timestamps.Length

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i

    $wsRaw.Cells.Item($row, 1).Value = 23
    $wsRaw.Cells.Item($row, 2).Value = 23

    # Column C holds the timestamp as text (cell is formatted as Text "@").
    # Make sure new rows (34-38) pick up the same Text number format the
    # existing rows use, since the column's default style is a time format.
    $cCell = $wsRaw.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $timestamps[$i]

    $wsRaw.Cells.Item($row, 4).Value = 0
    $wsRaw.Cells.Item($row, 5).Value = $temperatures[$i]
}

# "Fila final" input and dependent COUNT formula
$wsRaw.Range("H2").Value = 20
$wsRaw.Range("H3").Formula = "=COUNT(E:E)-7"

# Restore the selections recorded in the workbook views
$wsRaw.Range("H4").Select()
$wsValid.Activate()
$wsValid.Range("I20").Select()
